# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the worker/period table (B15:J29) on "Hoja1":
#  - Row 16 / Row 17 swap worker identity (Victor Alfonso Garcia Castro now
#    leads the table) and row 16 gets new Valor Mora / Salario Basico values.
#  - Rows 18-29 (Victor Enrique Villar Alvarez) now list the mora periods in
#    ascending order (1709, 1711, 1712, 1801-1809) instead of descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data doc id / name / period / valor mora / salario basico for rows 16-29.
$rows = @(
    @{ Row = 16; Doc = "1047372573"; Nombre = "VICTOR ALFONSO GARCIA CASTRO"; Periodo = "1701"; Valor = 27578;  Salario = 781242 },
    @{ Row = 17; Doc = "85487449";   Nombre = "DAIRO JOSE DE AVILA PADILLA";  Periodo = "1709"; Valor = 40000;  Salario = 1000000 },
    @{ Row = 18; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1709"; Valor = 40000; Salario = 1000000 },
    @{ Row = 19; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1711"; Valor = 40000; Salario = 1000000 },
    @{ Row = 20; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1712"; Valor = 40000; Salario = 1000000 },
    @{ Row = 21; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1801"; Valor = 40000; Salario = 1000000 },
    @{ Row = 22; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1802"; Valor = 40000; Salario = 1000000 },
    @{ Row = 23; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1803"; Valor = 40000; Salario = 1000000 },
    @{ Row = 24; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1804"; Valor = 40000; Salario = 1000000 },
    @{ Row = 25; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1805"; Valor = 40000; Salario = 1000000 },
    @{ Row = 26; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1806"; Valor = 40000; Salario = 1000000 },
    @{ Row = 27; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1807"; Valor = 40000; Salario = 1000000 },
    @{ Row = 28; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1808"; Valor = 40000; Salario = 1000000 },
    @{ Row = 29; Doc = "73198175";   Nombre = "VICTOR ENRIQUE VILLAR ALVAREZ"; Periodo = "1809"; Valor = 40000; Salario = 1000000 }
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 3).Value = $rec.Doc        # C - N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $rec.Nombre     # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $rec.Periodo    # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $rec.Valor      # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $rec.Salario    # G - Salario Basico
}
